# Generate Report for Handoff
#
# A new handoff report was generated for the source file
# "647f3c9e-28b9-4936-9d95-c0fe12ff52e7.md" (the first data row on every
# sheet). This refreshes that file's "Latest Handoff" timestamps:
#   - Overview!D2          "Latest Handoff Date"      -> 2016-03-30 09:37:27
#   - zh-cn!E2              "Latest Handoff Datetime"  -> 2016-03-30 09:37:19
#   - de-de!E2              "Latest Handoff Datetime"  -> 2016-03-30 09:37:27
# All other rows/cells are untouched.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("D2").Value = "2016-03-30 09:37:27"
$wsZhCn.Range("E2").Value     = "2016-03-30 09:37:19"
$wsDeDe.Range("E2").Value     = "2016-03-30 09:37:27"
